$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")
Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
